$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 15 new rows at the top of the data block (rows 2-16), pushing
# the existing comparison rows down by 15 (they become rows 17-46).
$ws.Range("A2:G16").Insert()

# Copy the formatting (cell styles: bold/centered/bordered A & B columns,
# default style on C:G) from the first surviving data block (now at
# rows 17:19) onto the newly inserted rows.
$src = $ws.Range("A17:G19")
$src.Copy()
$ws.Range("A2:G16").PasteSpecial(-4122)

# CRD
$ws.Range("A2").Value = "CRD"
$ws.Range("B2").Value = 70
$ws.Range("C2").Value = 0.002417
$ws.Range("D2").Value = -335.485371190243
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = -335.485371190243
$ws.Range("G2").Value = 0
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = 299.8
$ws.Range("D3").Value = -340.015794666801
$ws.Range("E3").Value = 30.9128228681901
$ws.Range("F3").Value = -309.10297179861
$ws.Range("G3").Value = 2.22966252682521
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1013
$ws.Range("D4").Value = -423.662018470503
$ws.Range("E4").Value = 272.86192163356
$ws.Range("F4").Value = -150.800096836943
$ws.Range("G4").Value = 1.81115134335482
$ws.Range("A2:A4").Merge()

# LBLRTM
$ws.Range("A5").Value = "LBLRTM"
$ws.Range("B5").Value = 70
$ws.Range("C5").Value = 0.0024
$ws.Range("D5").Value = -335.176214169739
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = -335.176214169739
$ws.Range("G5").Value = 0
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 299.8
$ws.Range("D6").Value = -339.80530455316
$ws.Range("E6").Value = 31.1053271108906
$ws.Range("F6").Value = -308.699977442271
$ws.Range("G6").Value = 2.09811707597551
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 1013
$ws.Range("D7").Value = -423.532488072306
$ws.Range("E7").Value = 272.902015332899
$ws.Range("F7").Value = -150.630472739406
$ws.Range("G7").Value = 1.76647688371748
$ws.Range("A5:A7").Merge()

# RRTMG
$ws.Range("A8").Value = "RRTMG"
$ws.Range("B8").Value = 70
$ws.Range("C8").Value = 0.00242
$ws.Range("D8").Value = -334.2143
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = -334.2143
$ws.Range("G8").Value = 0
$ws.Range("B9").Value = 30
$ws.Range("C9").Value = 299.7
$ws.Range("D9").Value = -339.1301
$ws.Range("E9").Value = 31.1483
$ws.Range("F9").Value = -307.9818
$ws.Range("G9").Value = 2.08088
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1013
$ws.Range("D10").Value = -423.5981
$ws.Range("E10").Value = 275.0453
$ws.Range("F10").Value = -148.5528
$ws.Range("G10").Value = 1.87119
$ws.Range("A8:A10").Merge()

# CLIRAD
$ws.Range("A11").Value = "CLIRAD"
$ws.Range("B11").Value = 70
$ws.Range("C11").Value = 0.002417
$ws.Range("D11").Value = -333.074616
$ws.Range("E11").Value = 0.00124102174077354
$ws.Range("F11").Value = -333.073374978259
$ws.Range("G11").Value = 5.88597492136063
$ws.Range("B12").Value = 30
$ws.Range("C12").Value = 299.8
$ws.Range("D12").Value = -337.951636
$ws.Range("E12").Value = 31.2661438252
$ws.Range("F12").Value = -306.6854921748
$ws.Range("G12").Value = 2.20048041954
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 1013
$ws.Range("D13").Value = -423.786188
$ws.Range("E13").Value = 275.551229
$ws.Range("F13").Value = -148.234959
$ws.Range("G13").Value = 1.898682161
$ws.Range("A11:A13").Merge()

# CLIRAD 96
$ws.Range("A14").Value = "CLIRAD 96"
$ws.Range("B14").Value = 70
$ws.Range("C14").Value = 0.002417
$ws.Range("D14").Value = -335.867256
$ws.Range("E14").Value = 0.00314171168088
$ws.Range("F14").Value = -335.864114288319
$ws.Range("G14").Value = 1.54851658
$ws.Range("B15").Value = 30
$ws.Range("C15").Value = 299.8
$ws.Range("D15").Value = -340.504712
$ws.Range("E15").Value = 29.5850774644
$ws.Range("F15").Value = -310.9196345356
$ws.Range("G15").Value = 2.101990674265
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 1013
$ws.Range("D16").Value = -423.786164
$ws.Range("E16").Value = 268.203585
$ws.Range("F16").Value = -155.582579
$ws.Range("G16").Value = 1.899698346
$ws.Range("A14:A16").Merge()

$ws.Range("A1").Select()
